# Rename the inline picture "name" metadata (wp:docPr/@name) for the three
# logo images that live in the document's headers/footers:
#   - First-page footer & default footer: Pearson logo  image1.png -> image2.png
#   - First-page header:                  BTec logo     image2.jpg -> image1.jpg
#
# The physical files in the package keep their original names/relationship
# ids - only the cosmetic "name" label on each picture changes.
#
# Note: InlineShape.Name only round-trips as a write (the getter does not
# reflect the existing docPr/@name), so shapes are targeted positionally /
# by AlternativeText rather than by reading back the current Name.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1 (the "default" header/footer)
# wdHeaderFooterFirstPage = 2 (the "first page" header/footer)
$headerFirst   = $sec.Headers.Item(2)
$footerDefault = $sec.Footers.Item(1)
$footerFirst   = $sec.Footers.Item(2)

# Pearson logo, first-page footer (was id="3" / name="image1.png")
$footerFirst.Range.InlineShapes.Item(1).Name = "image2.png"

# Pearson logo, default footer (was id="2" / name="image1.png")
$footerDefault.Range.InlineShapes.Item(1).Name = "image2.png"

# BTec logo, first-page header (was id="1" / name="image2.jpg")
$headerFirst.Range.InlineShapes.Item(1).Name = "image1.jpg"

Write-Output "Renamed inline picture names in headers/footers."
